$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    34 = 142852
    35 = 146674
    36 = 150987
    37 = 154815
    38 = 161141
    39 = 167127
    40 = 174812
    41 = 172465
    42 = 176954
    43 = 181677
    44 = 191265
    45 = 199751
    46 = 205236
    47 = 212157
    48 = 222547
    49 = 233093
    50 = 233518
    51 = 242699
    52 = 246197
    53 = 245171
    54 = 247775
    55 = 250819
    56 = 251535
    57 = 252963
    58 = 292451
    59 = 296588
    60 = 298595
    61 = 299207
    62 = 300445
    63 = 297894
    64 = 294384
    65 = 288538
    66 = 283563
    67 = 280520
    68 = 277716
    69 = 275690
    70 = 275275
    71 = 276590
    72 = 276744
    73 = 276911
    74 = 277219
    75 = 275189
    76 = 274244
    77 = 274302
    78 = 274535
    79 = 274455
    80 = 275050
    81 = 276050
    82 = 277112
    83 = 280276
    84 = 282123
    85 = 284786
    86 = 286347
    87 = 288308
    88 = 290213
    89 = 291699
    90 = 293879
    91 = 298142
    92 = 301282
    93 = 302146
    94 = 302368
    95 = 306592
    96 = 309669
    97 = 312901
    98 = 317031
    99 = 323217
    100 = 328307
    101 = 332431
    102 = 337724
    103 = 352949
    104 = 363116
    105 = 376492
    106 = 397284
    107 = 408772
    108 = 416274
    109 = 419861
    110 = 424786
    111 = 431406
    112 = 439391
    113 = 443619
    114 = 446737
    115 = 453431
    116 = 456388
    117 = 456618
    118 = 460207
    119 = 466376
    120 = 471580
}

foreach ($row in $updates.Keys) {
    $ws.Range("K$row").Value = $updates[$row]
}
